$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.080.48"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "3.568.01"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'235.74"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").Value = "'656.37"
$ws.Range("E6").Value = "  +2.39%  "
$ws.Range("D7").Value = "'1.47"
$ws.Range("E7").Value = "  -2.33%  "
$ws.Range("D8").Value = "'0.398"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -3.09%  "
$ws.Range("D11").Value = "3.567.19"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").Value = "'42.35"
$ws.Range("E13").Value = "  -3.15%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "4.233.21"
$ws.Range("E15").Value = "  -1.59%  "
$ws.Range("D16").Value = "94.976.17"
$ws.Range("E16").Value = "  -2.10%  "
$ws.Range("E17").Value = "  -0.95%  "
$ws.Range("D18").Value = "3.581.65"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "'7.75"
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("D20").Value = "'12.66"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").Value = "'17.76"
$ws.Range("E21").Value = "  -3.29%  "
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "'508.54"
$ws.Range("E23").Value = "  -1.76%  "
$ws.Range("D24").Value = "'0.479"
$ws.Range("E24").Value = "  -4.75%  "
$ws.Range("D25").Value = "'6.80"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("D27").Value = "'95.09"
$ws.Range("E27").Value = "  -3.56%  "
$ws.Range("D28").Value = "'12.64"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "3.760.26"
$ws.Range("E29").Value = "  -1.26%  "
$ws.Range("D30").Value = "'3.03"
$ws.Range("E30").Value = "  -3.02%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "'11.49"
$ws.Range("E32").Value = "  -1.98%  "
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("E35").Value = "  -3.62%  "
$ws.Range("D36").Value = "'31.87"
$ws.Range("E36").Value = "  +4.35%  "
$ws.Range("E37").Value = "  +11.28%  "
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("D39").Value = "'8.49"
$ws.Range("E39").Value = "  +7.12%  "
$ws.Range("D40").Value = "'578.66"
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("E44").Value = "  +3.43%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'34.82"
$ws.Range("E45").Value = "  +32.90%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "'5.75"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("E47").Value = "  +2.64%  "
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "'0.0413"
$ws.Range("E49").Value = "  -5.25%  "
$ws.Range("D50").Value = "'3.57"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "'8.15"
$ws.Range("E51").Value = "  -1.11%  "
